# "Ares Condicionados" is the 3rd sheet (matches dimension A1:H3 -> A1:H5,
# D2/E2 turning from numbers into literal "0" text, and 2 new rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ares Condicionados")

# Helper: force a literal text value (avoids Excel auto-coercing numeric-
# looking strings to numbers, or True/False to Booleans), then strips the
# resulting quote-prefix formatting so no stray style is left behind.
function Set-LiteralText {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2: D2 (Temperatura) and E2 (Intensidade) become the text "0" instead
# of the numbers 14 / 65.
Set-LiteralText $ws.Range("D2") "0"
Set-LiteralText $ws.Range("E2") "0"

# Row 3 is unchanged.

# New row 4.
$ws.Range("A4").Value = "teste"
$ws.Range("B4").Value = "teste"
Set-LiteralText $ws.Range("C4") "False"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

# New row 5.
$ws.Range("A5").Value = "ji"
$ws.Range("B5").Value = "vai da corinthinas"
Set-LiteralText $ws.Range("C5") "False"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
